$d = $word.ActiveDocument

# --- Replace Figure S3 caption (old para 5) with the new Figure S3 text ---
$d.Paragraphs(5).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S3. Protected land conversion pressure relative to unmanaged land for (b-d) reference and (f-h) low carbon transition scenarios across the four protection cases, normalized by the respective CURRENT case. These values are derived from those in Figure 4 by dividing each future protection case by the respective CURRENT case. Values represent individual land types within individual land units. The horizontal line is the median, the box represents the interquartile range, and the whiskers represent 1.5 times the interquartile range. The outliers are not shown.</w:t></w:r></w:p>')

# --- Replace Figure S4 caption (old para 7) with the new Figure S4 text ---
$d.Paragraphs(7).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S4. Relationship between protected land conversion pressure and suitable, protected area relative to managed area at the global level, for the low carbon transition BIODIV scenario (2015-2100).</w:t></w:r></w:p>')

# --- Old para 9 (previously "Figure S5...") becomes new Figure S5 (old S3 text), split into 3 runs ---
$d.Paragraphs(9).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S</w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t>. Relationships between protected land conversion pressure and suitable, protected area relative to managed area at the regional level, for the low carbon transition BIODIV scenario (2015-2100).</w:t></w:r></w:p>')

# --- Old para 11 (previously "Figure S6...") becomes new Figure S6 (old S4 text), split into 3 runs ---
$d.Paragraphs(11).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S</w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t>. Relationship between protected land conversion pressure and suitable, protected area relative to managed area for all individual land units, for the low carbon transition BIODIV scenario (2015-2100).</w:t></w:r></w:p>')

# --- Old para 13 (previously "Figure S7...") becomes new Figure S7 (new GCAM text) ---
$d.Paragraphs(13).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S7. GCAM global land allocation for the low carbon transition scenarios with different land protections.</w:t></w:r></w:p>')

# --- Old para 15 (previously "Figure S8...") becomes new Figure S8 (old S5 text), split into 3 runs ---
$d.Paragraphs(15).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S</w:t></w:r><w:r><w:t>8</w:t></w:r><w:r><w:t>. Distributions of regional difference from CURRENT as the percent of CURRENT allocation change from initial, for select years, under low carbon transition. The horizontal line is the median, the box represents the interquartile range, the whiskers represent 1.5 times the interquartile range, and the dots are the remaining outliers.</w:t></w:r></w:p>')

# --- Old para 17 (previously "Figure S9...") becomes new Figure S9 (old S6 text), split into 3 runs ---
$d.Paragraphs(17).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S</w:t></w:r><w:r><w:t>9</w:t></w:r><w:r><w:t>. Bioenergy crop production by region and globally for the low carbon transition scenarios with different land protections.</w:t></w:r></w:p>')

# --- Old trailing empty para 18 expands into: blank, S10, blank, S11, blank, S12, blank ---
$d.Paragraphs(18).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S</w:t></w:r><w:r><w:t>10</w:t></w:r><w:r><w:t>. Bioenergy electricity generation by region and globally for the low carbon transition scenarios with different land protections.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S</w:t></w:r><w:r><w:t>11</w:t></w:r><w:r><w:t>. Bioenergy consumption by region and globally for the low carbon transition scenarios with different land protections.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Figure S</w:t></w:r><w:r><w:t>12</w:t></w:r><w:r><w:t>. Total energy consumption by region and globally for the low carbon transition scenarios with different land protections.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# --- Prepend the new title paragraph + blank separator before everything ---
$d.Range(0, 0).InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Supplemental figure captions</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

